# Auto commit at 2025-11-25 8:39:43.56
#
# Updates the "Metrics" sheet's daily figures (B2:B13). The "today" sheet
# pulls these via formulas (=Metrics!B2 etc.), so its dependent cells
# recalc automatically. Also replays the user's navigation: ends up with
# the "today" sheet active (cell C8 selected) instead of "Chargingdata"
# (which keeps its old E33 selection but loses the active-tab flag), and
# leaves the "Metrics" sheet's remembered selection at F12.

$wb = $excel.ActiveWorkbook

# --- Metrics: refreshed daily figures -------------------------------------
$wsMetrics = $wb.Worksheets.Item("Metrics")

$wsMetrics.Range("B2").Value  = 325236.77000000008
$wsMetrics.Range("B3").Value  = 286549.01
$wsMetrics.Range("B4").Value  = 100303.46
$wsMetrics.Range("B5").Value  = 13256
$wsMetrics.Range("B6").Value  = 5121482.5200000014
$wsMetrics.Range("B7").Value  = 4328625.6900000004
$wsMetrics.Range("B8").Value  = 1507263.2900000003
$wsMetrics.Range("B9").Value  = 199463
$wsMetrics.Range("B10").Value = 33586863.510000013
$wsMetrics.Range("B11").Value = 31603900.850000005
$wsMetrics.Range("B12").Value = 11788985.33
$wsMetrics.Range("B13").Value = 1297093

# Leave Metrics' remembered cursor position at F12 (it is not the final
# active sheet, so this only updates its stored <selection>, not the tab).
$wsMetrics.Activate()
$wsMetrics.Range("F12").Select()

# --- Navigate to the "today" sheet, landing on C8 -------------------------
# This both clears tabSelected/E33-selection-owner from "Chargingdata"
# (previously the active tab) and marks "today" as active with the new
# selection, matching the workbook's activeTab index move (2 -> 5).
$wsToday = $wb.Worksheets.Item("today")
$wsToday.Activate()
$wsToday.Range("C8").Select()
